$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Entradas ")
$ws2 = $wb.Worksheets.Item("Saídas")

# --- Sheet "Entradas " (sheet1): extend the sponsorship (Patrocínio/Ambev) rows
# Row 3 already holds the first Patrocínio/Ambev entry; clone its formatting
# down through row 10 for the new entries, then fill in the per-row values.
$ws1.Range("A3:H3").Copy()
$ws1.Range("A4:H10").PasteSpecial(-4122)

$ws1.Cells.Item(4, 1).Value = "Patrocínio"
$ws1.Cells.Item(4, 2).Value = 44565
$ws1.Cells.Item(4, 3).Value = 5501
$ws1.Cells.Item(4, 4).Value = "Ambev"
$ws1.Cells.Item(4, 5).Value = 2
$ws1.Cells.Item(4, 6).Value = "-"
$ws1.Cells.Item(4, 7).Value = "-"
$ws1.Cells.Item(4, 8).Value = "-"

$ws1.Cells.Item(5, 1).Value = "Patrocínio"
$ws1.Cells.Item(5, 2).Value = 44778
$ws1.Cells.Item(5, 3).Value = 5502
$ws1.Cells.Item(5, 4).Value = "Ambev"
$ws1.Cells.Item(5, 5).Value = 3
$ws1.Cells.Item(5, 6).Value = "-"
$ws1.Cells.Item(5, 7).Value = "-"
$ws1.Cells.Item(5, 8).Value = "-"

$ws1.Cells.Item(6, 1).Value = "Patrocínio"
$ws1.Cells.Item(6, 2).Value = 44779
$ws1.Cells.Item(6, 3).Value = 5503
$ws1.Cells.Item(6, 4).Value = "Ambev"
$ws1.Cells.Item(6, 5).Value = 4
$ws1.Cells.Item(6, 6).Value = "-"
$ws1.Cells.Item(6, 7).Value = "-"
$ws1.Cells.Item(6, 8).Value = "-"

$ws1.Cells.Item(7, 1).Value = "Patrocínio"
$ws1.Cells.Item(7, 2).Value = 44780
$ws1.Cells.Item(7, 3).Value = 5504
$ws1.Cells.Item(7, 4).Value = "Ambev"
$ws1.Cells.Item(7, 5).Value = 5
$ws1.Cells.Item(7, 6).Value = "-"
$ws1.Cells.Item(7, 7).Value = "-"
$ws1.Cells.Item(7, 8).Value = "-"

$ws1.Cells.Item(8, 1).Value = "Patrocínio"
$ws1.Cells.Item(8, 2).Value = 44781
$ws1.Cells.Item(8, 3).Value = 5505
$ws1.Cells.Item(8, 4).Value = "Ambev"
$ws1.Cells.Item(8, 5).Value = 6
$ws1.Cells.Item(8, 6).Value = "-"
$ws1.Cells.Item(8, 7).Value = "-"
$ws1.Cells.Item(8, 8).Value = "-"

$ws1.Cells.Item(9, 1).Value = "Patrocínio"
$ws1.Cells.Item(9, 2).Value = 44570
$ws1.Cells.Item(9, 3).Value = 5506
$ws1.Cells.Item(9, 4).Value = "Ambev"
$ws1.Cells.Item(9, 5).Value = 7
$ws1.Cells.Item(9, 6).Value = "-"
$ws1.Cells.Item(9, 7).Value = "-"
$ws1.Cells.Item(9, 8).Value = "-"

$ws1.Cells.Item(10, 1).Value = "Patrocínio"
$ws1.Cells.Item(10, 2).Value = 44571
$ws1.Cells.Item(10, 3).Value = 5507
$ws1.Cells.Item(10, 4).Value = "Ambev"
$ws1.Cells.Item(10, 5).Value = 8
$ws1.Cells.Item(10, 6).Value = "-"
$ws1.Cells.Item(10, 7).Value = "-"
$ws1.Cells.Item(10, 8).Value = "-"

# --- Sheet "Saídas" (sheet2): extend the "Fornecedor de destilados" payments
# Row 3's "02/012020" text placeholder becomes a real date, and rows 4-7 are
# new payment entries cloned from row 3's formatting.
$ws2.Range("A3:E3").Copy()
$ws2.Range("A4:E7").PasteSpecial(-4122)

$ws2.Cells.Item(3, 2).Value = 44563

$ws2.Cells.Item(4, 1).Value = "Pagamento"
$ws2.Cells.Item(4, 2).Value = 44564
$ws2.Cells.Item(4, 3).Value = 1564.52
$ws2.Cells.Item(4, 4).Value = "João Flávio"
$ws2.Cells.Item(4, 5).Value = "Fornecedor de destilados"

$ws2.Cells.Item(5, 1).Value = "Pagamento"
$ws2.Cells.Item(5, 2).Value = 44565
$ws2.Cells.Item(5, 3).Value = 1565.52
$ws2.Cells.Item(5, 4).Value = "João Flávio"
$ws2.Cells.Item(5, 5).Value = "Fornecedor de destilados"

$ws2.Cells.Item(6, 1).Value = "Pagamento"
$ws2.Cells.Item(6, 2).Value = 44778
$ws2.Cells.Item(6, 3).Value = 1566.52
$ws2.Cells.Item(6, 4).Value = "João Flávio"
$ws2.Cells.Item(6, 5).Value = "Fornecedor de destilados"

$ws2.Cells.Item(7, 1).Value = "Pagamento"
$ws2.Cells.Item(7, 2).Value = 44779
$ws2.Cells.Item(7, 3).Value = 1567.52
$ws2.Cells.Item(7, 4).Value = "João Flávio"
$ws2.Cells.Item(7, 5).Value = "Fornecedor de destilados"

# --- Selections: match the saved cursor position on each sheet; select
# sheet1's cell last-but-one so sheet2 (tab-selected in the source file)
# ends up the active tab after this script runs.
$ws1.Range("D13").Select()
$ws2.Range("B3").Select()
